# "Update countries & provincias Spain"
# Refresh the COVID-19 country table on sheet "Pais":
#  - bump the "last updated" timestamp
#  - push newer case/death/recovered counts into the countries whose
#    figures moved since the previous pull
#  - two countries (Ucrania / Emiratos Arabes Unidos and
#    Haiti / Estado de Palestina) leapfrogged their neighbour in the
#    total-cases ranking, so those row pairs swap place; Islas Malvinas
#    and Groenlandia are tied on totals and also swap order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- timestamp ---------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 15 de Julio de 2020 a las 09:01"

# --- plain data refreshes (country keeps its row) -----------------------
# Row 4: Estados Unidos
$ws.Range("B4").Value = 3545254
$ws.Range("C4").Value = 177
$ws.Range("D4").Value = 1600321
$ws.Range("E4").Value = 1805788
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 139145

# Row 6: India
$ws.Range("B6").Value = 937844
$ws.Range("C6").Value = 357
$ws.Range("D6").Value = 593178
$ws.Range("E6").Value = 320339
$ws.Range("G6").Value = 12
$ws.Range("H6").Value = 24327

# Row 48: Afganistan
$ws.Range("B48").Value = 34994
$ws.Range("C48").Value = 254
$ws.Range("D48").Value = 22456
$ws.Range("E48").Value = 11444
$ws.Range("G48").Value = 46
$ws.Range("H48").Value = 1094

# Row 53: Armenia
$ws.Range("B53").Value = 33005
$ws.Range("C53").Value = 515
$ws.Range("D53").Value = 21348
$ws.Range("E53").Value = 11065
$ws.Range("G53").Value = 11
$ws.Range("H53").Value = 592

# Row 55: Honduras
$ws.Range("D55").Value = 3287
$ws.Range("E55").Value = 25012

# Row 76: El Salvador
$ws.Range("D76").Value = 5947
$ws.Range("E76").Value = 4078

# Row 99: Hungria
$ws.Range("B99").Value = 4263
$ws.Range("C99").Value = 5
$ws.Range("D99").Value = 3126
$ws.Range("E99").Value = 542

# Row 116: Sudan del Sur
$ws.Range("B116").Value = 2153
$ws.Range("C116").Value = 5
$ws.Range("D116").Value = 1175
$ws.Range("E116").Value = 937

# Row 124: Cabo Verde
$ws.Range("B124").Value = 1780
$ws.Range("C124").Value = 58
$ws.Range("D124").Value = 850
$ws.Range("E124").Value = 911

# Row 144: Georgia
$ws.Range("B144").Value = 1003
$ws.Range("C144").Value = 4
$ws.Range("D144").Value = 873
$ws.Range("E144").Value = 115

# Row 185: Seychelles
$ws.Range("D185").Value = 27
$ws.Range("E185").Value = 73

# --- ranking swaps -------------------------------------------------------
# Ucrania's refreshed total (55607) now exceeds Emiratos Arabes Unidos'
# previous total (55573), so Ucrania takes row 38 with its new figures and
# Emiratos Arabes Unidos drops to row 39 carrying the figures it already had.
$ws.Range("A38").Value = "Ucrania"
$ws.Range("B38").Value = 55607
$ws.Range("C38").Value = 836
$ws.Range("D38").Value = 28131
$ws.Range("E38").Value = 26049
$ws.Range("G38").Value = 15
$ws.Range("H38").Value = 1427

$ws.Range("A39").Value = "Emiratos Arabes Unidos"
$ws.Range("B39").Value = 55573
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 46025
$ws.Range("E39").Value = 9213
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 335

# Haiti's refreshed total (6831) now exceeds Estado de Palestina's previous
# total (6764), so Haiti takes row 88 and Estado de Palestina drops to row 89.
$ws.Range("A88").Value = "Haiti"
$ws.Range("B88").Value = 6831
$ws.Range("C88").Value = 104
$ws.Range("D88").Value = 3283
$ws.Range("E88").Value = 3405
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 143

$ws.Range("A89").Value = "Estado de Palestina"
$ws.Range("B89").Value = 6764
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 1084
$ws.Range("E89").Value = 5636
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 44

# Islas Malvinas and Groenlandia are tied on every figure (13/0/13/0/0/0/0);
# they simply swap display order.
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"
